# Zeiterfassung.xlsx - "Projektdokumentation, Anpassung Phasen, Zeiterfassung"
#
# 1. Add three new time-tracking entries (rows 4-6) on the "Zeiterfassung" sheet.
# 2. Make "Zeiterfassung" the active sheet/tab and set the new selection on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# --- New rows of tracked time -------------------------------------------------
$ws.Range("B4").Value = "Brainstorming, telefonat mit Mejdin"
$ws.Range("C4").Value = "Tobias Lanz"
$ws.Range("D4").Value = 42245
$ws.Range("F4").Value = 0.75

$ws.Range("B5").Value = "Mindmap Inputs"
$ws.Range("C5").Value = "Tobias Lanz"
$ws.Range("D5").Value = 42246
$ws.Range("F5").Value = 1.5

$ws.Range("B6").Value = "Initiale Projektdokumentation"
$ws.Range("C6").Value = "Tobias Lanz"
$ws.Range("D6").Value = 42246
$ws.Range("F6").Value = 2.5

# --- Make "Zeiterfassung" the active sheet / update its selection -------------
$ws.Activate() | Out-Null
$ws.Range("D9").Select() | Out-Null
